$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -17
$ws.Range("F3").Value = -4
$ws.Range("F4").Value = -5
$ws.Range("F5").Value = -9
$ws.Range("F6").Value = -2
$ws.Range("F8").Value = -9
$ws.Range("F11").Value = -3
$ws.Range("F12").Value = 0
$ws.Range("F16").Value = -1
$ws.Range("F18").Value = -5
$ws.Range("F20").Value = 11
$ws.Range("F23").Value = -2
$ws.Range("F27").Value = 3
